# Apply updated cryptocurrency market data to the sheet.
# Only the cells whose content actually changed are listed below,
# keyed by row number with the column letters that were updated.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row=2; D="24.861.92"; E="  -4.20%  " },
    @{ Row=3; D="1.633.05"; E="  -6.48%  " },
    @{ Row=4; D="0.9981"; E="  -0.16%  " },
    @{ Row=5; D="234.43"; E="  -5.83%  " },
    @{ Row=6; E="  +0.04%  " },
    @{ Row=7; D="0.4718"; E="  -6.65%  " },
    @{ Row=8; B="Cardano"; C="https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"; D="0.2556"; E="  -7.02%  " },
    @{ Row=9; B="Dogecoin"; C="https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"; D="0.06079"; E="  -1.73%  " },
    @{ Row=10; B="TRON"; C="https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"; D="0.06943"; E="  -4.41%  " },
    @{ Row=11; B="WrappedEther"; C="https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"; D="1.638.89"; E="  -6.09%  " },
    @{ Row=12; B="Solana"; C="https://coinranking.com/coin/zNZHO_Sjf+solana-sol"; D="14.57"; E="  -4.14%  " },
    @{ Row=13; B="Polygon"; C="https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"; D="0.6010"; E="  -8.19%  " },
    @{ Row=14; B="Polkadot"; C="https://coinranking.com/coin/25W7FG7om+polkadot-dot"; D="4.312"; E="  -7.32%  " },
    @{ Row=15; B="Litecoin"; C="https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"; D="72.97"; E="  -6.10%  " },
    @{ Row=16; B="Dai"; C="https://coinranking.com/coin/MoTuySvg7+dai-dai"; D="1.000"; E="  +0.04%  " },
    @{ Row=17; B="BinanceUSD"; C="https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"; D="0.9986"; E="  -0.06%  " },
    @{ Row=18; B="WrappedBTC"; C="https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"; D="24.862.08"; E="  -4.29%  " },
    @{ Row=19; B="ShibaInu"; C="https://coinranking.com/coin/xz24e0BjL+shibainu-shib"; D="0.000006556"; E="  -4.18%  " },
    @{ Row=20; B="Avalanche"; C="https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"; D="11.11"; E="  -6.16%  " },
    @{ Row=21; B="WrappedliquidstakedEther2.0"; C="https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"; D="1.850.19"; E="  -5.99%  " },
    @{ Row=22; B="Uniswap"; C="https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"; D="4.340"; E="  -2.55%  " },
    @{ Row=23; B="Cosmos"; C="https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"; D="8.537"; E="  -2.20%  " },
    @{ Row=24; B="Chainlink"; C="https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"; D="5.208"; E="  -3.51%  " },
    @{ Row=25; B="Monero"; C="https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"; D="133.01"; E="  -2.76%  " },
    @{ Row=26; B="EthereumClassic"; C="https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"; D="14.73"; E="  -3.42%  " },
    @{ Row=27; B="Toncoin"; C="https://coinranking.com/coin/67YlI0K1b+toncoin-ton"; D="1.374"; E="  -8.79%  " },
    @{ Row=28; B="BitcoinCash"; C="https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"; D="103.15"; E="  -2.37%  " },
    @{ Row=29; B="LidoDAOToken"; C="https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"; D="1.620"; E="  -9.21%  " },
    @{ Row=30; B="InternetComputer(DFINITY)"; C="https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"; D="3.785"; E="  -2.28%  " },
    @{ Row=31; B="Stellar"; C="https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"; D="0.07702"; E="  -6.22%  " },
    @{ Row=32; B="Filecoin"; C="https://coinranking.com/coin/ymQub4fuB+filecoin-fil"; D="3.519"; E="  -3.54%  " },
    @{ Row=33; B="Frax"; C="https://coinranking.com/coin/KfWtaeV1W+frax-frax"; D="0.9992"; E="  +0.02%  " },
    @{ Row=34; B="Hedera"; C="https://coinranking.com/coin/jad286TjB+hedera-hbar"; D="0.04270"; E="  -8.75%  " },
    @{ Row=35; B="HuobiToken"; C="https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"; D="2.581"; E="  -2.80%  " },
    @{ Row=36; B="ARBITRUM"; C="https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"; D="0.9188"; E="  -7.73%  " },
    @{ Row=37; B="ImmutableX"; C="https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; D="0.5766"; E="  -6.89%  " },
    @{ Row=38; B="MXToken"; C="https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"; D="2.539"; E="  -7.92%  " },
    @{ Row=39; B="VeChain"; C="https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; D="0.01536"; E="  -4.77%  " },
    @{ Row=40; B="PaxDollar"; C="https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"; D="0.9989"; E="  -0.07%  " },
    @{ Row=41; B="TrustWalletToken"; C="https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"; D="0.8060"; E="  +5.96%  " },
    @{ Row=42; B="Quant"; C="https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"; D="97.23"; E="  -3.39%  " },
    @{ Row=43; B="RenderToken"; C="https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"; D="1.758"; E="  -8.91%  " },
    @{ Row=44; B="TheSandbox"; C="https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"; D="0.3679"; E="  -6.50%  " },
    @{ Row=45; B="FraxShare"; C="https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"; D="4.695"; E="  -6.26%  " },
    @{ Row=46; B="Cronos"; C="https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"; D="0.05196"; E="  -1.61%  " },
    @{ Row=47; B="Algorand"; C="https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"; D="0.1088"; E="  -5.66%  " },
    @{ Row=48; B="Aptos"; C="https://coinranking.com/coin/HGYj5JCv5+aptos-apt"; D="5.998"; E="  -5.14%  " },
    @{ Row=49; B="Elrond"; C="https://coinranking.com/coin/omwkOTglq+elrond-egld"; D="29.31"; E="  -4.50%  " },
    @{ Row=50; B="TrueUSD"; C="https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd"; D="1.000"; E="  -0.20%  " },
    @{ Row=51; B="USDD"; C="https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"; D="0.9958"; E="  -0.27%  " }
)

foreach ($item in $updates) {
    $r = $item.Row
    if ($item.ContainsKey("B")) { $ws.Range("B$r").Value = $item.B }
    if ($item.ContainsKey("C")) { $ws.Range("C$r").Value = $item.C }
    if ($item.ContainsKey("D")) {
        # Column D (Price) can contain digit-only text (e.g. "0.9981").
        # Format the cell as Text first so Excel keeps the exact string
        # instead of silently converting it to a numeric value.
        $ws.Range("D$r").NumberFormat = "@"
        $ws.Range("D$r").Value = $item.D
    }
    if ($item.ContainsKey("E")) { $ws.Range("E$r").Value = $item.E }
}